# Update cryptocurrency price/volume figures (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'93.657.47"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.09%  "
$ws.Range("D3").Value = "'3.045.94"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.74%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'232.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.70%  "
$ws.Range("D6").Value = "'602.74"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.42%  "
$ws.Range("D7").Value = "'1.08"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.70%  "
$ws.Range("D8").Value = "'0.370"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -9.70%  "
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("D10").Value = "'0.795"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.54%  "
$ws.Range("D11").Value = "'3.039.49"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.88%  "
$ws.Range("E12").Value = "  -4.14%  "
$ws.Range("D13").Value = "'93.230.91"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.10%  "
$ws.Range("D14").Value = "'0.0000236"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -8.05%  "
$ws.Range("D15").Value = "'33.01"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.52%  "
$ws.Range("D16").Value = "'5.23"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.82%  "
$ws.Range("D17").Value = "'3.608.85"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.13%  "
$ws.Range("D18").Value = "'3.049.85"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.72%  "
$ws.Range("D19").Value = "'3.45"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -9.19%  "
$ws.Range("D20").Value = "'14.16"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.27%  "
$ws.Range("D21").Value = "'5.61"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.45%  "
$ws.Range("D22").Value = "'432.44"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.28%  "
$ws.Range("D23").Value = "'8.66"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -8.28%  "
$ws.Range("E24").Value = "  -11.68%  "
$ws.Range("D25").Value = "'8.34"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.42%  "
$ws.Range("D26").Value = "'5.42"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -7.46%  "
$ws.Range("D27").Value = "'83.79"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.73%  "
$ws.Range("D28").Value = "'11.55"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.95%  "
$ws.Range("D29").Value = "'3.196.53"
$ws.Range("D29").Style = "Normal"
$ws.Range("E30").Value = "  +0.03%  "
$ws.Range("D31").Value = "'0.243"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.91%  "
$ws.Range("D32").Value = "'1.12"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +11.99%  "
$ws.Range("D33").Value = "'0.174"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.04%  "
$ws.Range("D34").Value = "'0.120"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -12.96%  "
$ws.Range("D35").Value = "'8.89"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.75%  "
$ws.Range("D36").Value = "'7.54"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -7.30%  "
$ws.Range("D37").Value = "'0.151"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.71%  "
$ws.Range("D38").Value = "'24.98"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.51%  "
$ws.Range("E39").Value = "  -2.73%  "
$ws.Range("D40").Value = "'23.96"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.71%  "
$ws.Range("D41").Value = "'3.77"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.87%  "
$ws.Range("D42").Value = "'0.426"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.60%  "
$ws.Range("D43").Value = "'455.31"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.31%  "
$ws.Range("E44").Value = "  -6.89%  "
$ws.Range("D46").Value = "'3.07"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -12.10%  "
$ws.Range("D47").Value = "'160.80"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.54%  "
$ws.Range("D48").Value = "'1.79"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -8.05%  "
$ws.Range("D49").Value = "'0.653"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.44%  "
$ws.Range("D50").Value = "'43.57"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.07%  "
$ws.Range("E51").Value = "  +0.01%  "
